{"js": "// The commit updates the \"Repository Usage Agreement Report\" consent\n// paragraph for the Third Call:\n//   1. \"each of the members of the C1.037\" -> \"the members of the C1.037\"\n//   2. \"... groups shows their agreement ...\" -> \"... groups show their agreement ...\"\n//      (subject/verb agreement fix to match the now-plural subject)\n//\n// Both edits live in the same paragraph, each on a distinct, unique\n// substring, so a simple body-level search & replace is sufficient and\n// keeps the surrounding run formatting (sz/szCs) untouched.\n\nconst body = context.document.body;\n\n// --- Edit 1: \"each of the members of\" -> \"the members of\" ---\nconst oldLeadIn = \"By signing this document, each of the members of the C1.037\";\nconst newLeadIn = \"By signing this document, the members of the C1.037\";\n\nconst leadInResults = body.search(oldLeadIn, { matchCase: true, matchWholeWord: false });\nleadInResults.load(\"items\");\nawait context.sync();\n\nif (leadInResults.items.length > 0) {\n  leadInResults.items[0].insertText(newLeadIn, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Edit 2: \"shows their agreement ... subject.\" -> \"show their agreement ... subject.\" ---\nconst oldTail =\n  \" shows their agreement to allow the work previously carried out during the regular term of the Design and Testing 2 course to be used for submission in the second and third calls. This consent confirms that they have no objections to the reuse of the project and understand their work done may be used by their teammates in future calls of the subject.\";\nconst newTail =\n  \" show their agreement to allow the work previously carried out during the regular term of the Design and Testing 2 course to be used for submission in the second and third calls. This consent confirms that they have no objections to the reuse of the project and understand their work done may be used by their teammates in future calls of the subject.\";\n\nconst tailResults = body.search(oldTail, { matchCase: true, matchWholeWord: false });\ntailResults.load(\"items\");\nawait context.sync();\n\nif (tailResults.items.length > 0) {\n  tailResults.items[0].insertText(newTail, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# \"readme modified and reports updated to Third Call\"\n#\n# In the \"Repository Usage Agreement Report\", the consent paragraph is\n# reworded for subject/verb agreement now that it addresses the group\n# members collectively:\n#   1. \"each of the members of the C1.037\" -> \"the members of the C1.037\"\n#   2. \"... groups shows their agreement ...\" -> \"... groups show their agreement ...\"\n#\n# Both changes live in the same paragraph. Plain Range.Text assignment in\n# this host coalesces every run in the touched paragraph (they all share\n# identical rPr), which would collapse the paragraph's existing run\n# boundaries beyond what actually changed. To keep the untouched runs\n# (\" and C2.037\", \" group\", \"s\") intact and land exactly on the wording\n# above, the paragraph's text (everything except its trailing paragraph\n# mark, so w:pPr / paragraph identity are left alone) is replaced via\n# Range.InsertXML with the equivalent WordprocessingML runs.\n\n$d = $word.ActiveDocument\n\n# Locate the consent paragraph by its unique lead-in text, then expand\n# the found range out to the whole enclosing paragraph.\n$found = $d.Content\n$found.Find.ClearFormatting()\n$found.Find.MatchCase = $true\n$ok = $found.Find.Execute(\"By signing this document, each of the members of the C1.037\")\nif (-not $ok) {\n    throw \"Could not locate the consent paragraph to update\"\n}\n$found.Expand(4) | Out-Null  # wdParagraph\n\n# Range over the paragraph's text only (End - 1 excludes the trailing\n# paragraph mark).\n$textRange = $d.Range($found.Start, $found.End - 1)\n\n$newRunsXml = @'\n<w:r w:rsidRPr=\"008C4D8B\"><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>By signing this document, the members of the C1.037</w:t></w:r><w:r w:rsidR=\"00D32969\"><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> and C2.037</w:t></w:r><w:r w:rsidRPr=\"008C4D8B\"><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> group</w:t></w:r><w:r w:rsidR=\"00D62736\"><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>s</w:t></w:r><w:r w:rsidRPr=\"008C4D8B\"><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> show their agreement to allow the work previously carried out during the regular term of the Design and Testing 2 course to be used for submission in the second and third calls. This consent confirms that they have no objections to the reuse of the project and understand their work done may be used by their teammates in future calls of the subject.</w:t></w:r>\n'@\n\n$pkg = @\"\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<?mso-application progid=\"Word.Document\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body><w:p>$newRunsXml</w:p></w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n\"@\n\n$textRange.InsertXML($pkg)\n"}
